# Add 6 new bet rows (rows 10-15) to the "bets" sheet.
#
# Each row is produced by duplicating the row directly above it (so the
# D/F/M running-total & percentage formulas -- and their cell styles --
# carry down exactly the way a manual "copy row / insert copied cells"
# in Excel would), then the input cells (id, date, W/L, stake, esport,
# teams, bet type, W, L) are overwritten with the new bet's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# Row data in sheet order (row 10 .. row 15).
$newRows = @(
    @{ Row=10; A=9;  B=45213; C=0; E=-5000; G="ESPORTS"; H="WORLDS 2023";               I="LOUD";  J="GANA 1 MAPA EN LA SERIE"; K=0; L=1 },
    @{ Row=11; A=10; B=45213; C=1; E=1140;  G="ESPORTS"; H="HALO WORLDS 2023";          I="OPTIC"; J="GANA SERIE";             K=1; L=0 },
    @{ Row=12; A=11; B=45213; C=1; E=690;   G="ESPORTS"; H="VALORANT CHAMPIONS LATAM";  I="KRU";   J="GANA SERIE";             K=1; L=0 },
    @{ Row=13; A=12; B=45214; C=1; E=64;    G="ESPORTS"; H="DOTA 2 THE INTERNATIONAL";  I="LGD";   J="GANA SERIE";             K=1; L=0 },
    @{ Row=14; A=13; B=45214; C=1; E=605;   G="ESPORTS"; H="DOTA 2 THE INTERNATIONAL";  I="LGD";   J="GANA SERIE";             K=1; L=0 },
    @{ Row=15; A=14; B=45214; C=1; E=3872;  G="ESPORTS"; H="DOTA 2 THE INTERNATIONAL";  I="LGD";   J="GANA SERIE";             K=1; L=0 }
)

# Pass 1: create the rows top-down by duplicating the row immediately above
# each new one. Doing this forwards (instead of always inserting at row 10)
# keeps every copied formula's relative reference pointing at the correct
# "row above" once all six exist.
foreach ($r in $newRows) {
    $prev = $r.Row - 1
    $ws.Rows($prev).Copy()
    $ws.Rows($r.Row).Insert()
}

# Pass 2: fill in the non-team input columns, forwards.
foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("C" + $r.Row).Value = $r.C
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("J" + $r.Row).Value = $r.J
    $ws.Range("K" + $r.Row).Value = $r.K
    $ws.Range("L" + $r.Row).Value = $r.L
}

# Pass 3: fill in the team-name columns (H, I) bottom-up -- this is the
# order the bets were actually typed in (most recent match entered
# first), which is also the order their names first appear in the
# workbook's shared-string table.
for ($i = $newRows.Count - 1; $i -ge 0; $i--) {
    $r = $newRows[$i]
    $ws.Range("H" + $r.Row).Value = $r.H
    $ws.Range("I" + $r.Row).Value = $r.I
}

# Move the active selection to reflect where the user ended up after adding
# the rows.
$ws.Range("F11").Select() | Out-Null
